$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("master-location")

# New master-location rows added for the "22 April Changes" update
# (Postal Code hierarchy entries under parent BNMR, in eng/fra/ara).
$newRows = @(
    @(110, 10110, 10110, 5, "Postal Code",    "BNMR", "eng"),
    @(111, 10111, 10111, 5, "Postal Code",    "BNMR", "eng"),
    @(112, 10113, 10113, 5, "Postal Code",    "BNMR", "eng"),
    @(113, 10114, 10114, 5, "Postal Code",    "BNMR", "eng"),
    @(114, 10111, 10111, 5, "code postal",    "BNMR", "fra"),
    @(115, 10110, 10110, 5, "code postal",    "BNMR", "fra"),
    @(116, 10113, 10113, 5, "code postal",    "BNMR", "fra"),
    @(117, 10114, 10114, 5, "code postal",    "BNMR", "fra"),
    @(118, 10111, 10111, 5, "الرمز البريدي",  "BNMR", "ara"),
    @(119, 10110, 10110, 5, "الرمز البريدي",  "BNMR", "ara")
)

foreach ($row in $newRows) {
    $r = $row[0]
    $ws.Cells.Item($r, 1).Value = $row[1]
    $ws.Cells.Item($r, 2).Value = $row[2]
    $ws.Cells.Item($r, 3).Value = $row[3]
    $ws.Cells.Item($r, 4).Value = $row[4]
    $ws.Cells.Item($r, 5).Value = $row[5]
    $ws.Cells.Item($r, 6).Value = $row[6]
    $ws.Cells.Item($r, 7).Value = $true
    $ws.Cells.Item($r, 8).Value = "superadmin"
    $ws.Cells.Item($r, 9).Value = "now()"
}
